$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 57 data. Several columns hold values that *look* numeric/date-like
# ("-499", "7/3/2025", "4", "1") but must be stored as TEXT (matching the
# inlineStr cells in the source diff), so force text via NumberFormat "@"
# before assigning, then clear the formatting override so no extra style
# index is left attached to the cell.
$ws.Range("A57").NumberFormat = "@"
$ws.Range("A57").Value = "-499"
$ws.Range("A57").ClearFormats()

$ws.Range("B57").NumberFormat = "@"
$ws.Range("B57").Value = "7/3/2025"
$ws.Range("B57").ClearFormats()

$ws.Range("C57").Value = "La Rioja 1770"

$ws.Range("D57").NumberFormat = "@"
$ws.Range("D57").Value = "4"
$ws.Range("D57").ClearFormats()

# E57 is present but empty in the source data.
$ws.Range("E57").NumberFormat = "@"
$ws.Range("E57").Value = ""
$ws.Range("E57").ClearFormats()

$ws.Range("F57").Value = "PEBCOM"
$ws.Range("G57").Value = "Pendiente"
$ws.Range("H57").Value = "picada"

$ws.Range("I57").NumberFormat = "@"
$ws.Range("I57").Value = "1"
$ws.Range("I57").ClearFormats()

$ws.Range("J57").Value = "Cambio"
$ws.Range("K57").Value = "Sin equipos"
$ws.Range("L57").Value = "Pasante"

$ws.Range("M57").Value = -58.406225
$ws.Range("N57").Value = -34.631314

$ws.Range("O57").Value = "San Telmo"
$ws.Range("P57").Value = "Capital Sur"
